# Generate Report for Handoff
#
# Updates the localization-status workbook:
#  - refreshes the existing "bf2ba36c..." handoff (row 2) to the new
#    "12dce014..." handoff markdown file, with refreshed handoff/target
#    timestamps
#  - appends two brand-new dependency rows (row 3 & row 4) for the two
#    .png assets that came along with this handoff, on all three sheets
#    (Overview, zh-cn, de-de)

$wb = $excel.ActiveWorkbook

$Missing = [System.Type]::Missing

$mdName  = "12dce014-7309-4ba6-8ff7-4e9e3fa1cb91.md"
$png1    = "26b0abf3-2f3e-4d5f-bb2f-97db13a43631.png"
$png2    = "6b7b7507-b4f4-43bf-89a2-0340226c9aa3.png"

$xlfZh   = "12dce014-7309-4ba6-8ff7-4e9e3fa1cb91.06e4b380a67ce395598a12746074389663b27a31.zh-cn.xlf"
$xlfDe   = "12dce014-7309-4ba6-8ff7-4e9e3fa1cb91.06e4b380a67ce395598a12746074389663b27a31.de-de.xlf"
$png1TargetZh = "23844c198639e6e12d86c4851b53f29446de8bcd.png"
$png2TargetZh = "f69815e00114aecbeb4069bd77f9916f46bf8710.png"

$handoffDate = "2016-08-13 19:08:35"
$zhDate      = "2016-03-13 19:08:32"
$deDate      = "2016-03-13 19:08:35"
$epoch       = "0001-01-01 00:00:00"

$readyStatus = "Ready for handoff"
$includeReason = "Include"
$isDependencyReason = "IsDependency"
$dependencyFrom = "e2e\" + $mdName

$mdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/724aae5eb971fb21d781924a5fb3dc7df0e1c19d/e2e/" + $mdName
$png1Url = "https://github.com/OpenLocalizationTest/oltest/blob/724aae5eb971fb21d781924a5fb3dc7df0e1c19d/e2e/" + $png1
$png2Url = "https://github.com/OpenLocalizationTest/oltest/blob/724aae5eb971fb21d781924a5fb3dc7df0e1c19d/e2e/" + $png2

$xlfZhUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/979a94049397f43f59ef6ce98aa5333abf56a735/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $xlfZh
$png1ZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/979a94049397f43f59ef6ce98aa5333abf56a735/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $png1TargetZh
$png2ZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/979a94049397f43f59ef6ce98aa5333abf56a735/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $png2TargetZh

$xlfDeUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f889a7fecbc757a73dbd23347126f72bb237d664/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $xlfDe
$png1DeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f889a7fecbc757a73dbd23347126f72bb237d664/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $png1TargetZh
$png2DeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f889a7fecbc757a73dbd23347126f72bb237d664/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $png2TargetZh

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 2: refresh the existing handoff file name + handoff date text.
# (Only the cell text is changed - the existing hyperlink on A2 is left
# alone so it keeps its single relationship/style instead of being
# duplicated.)
$wsOverview.Range("A2").Value2 = $mdName
$wsOverview.Range("B2").Value2 = $readyStatus
$wsOverview.Range("C2").Value2 = $readyStatus
$wsOverview.Range("D2").Value2 = $handoffDate

# Row 3: first new dependency (png)
$wsOverview.Range("B3").Value2 = $readyStatus
$wsOverview.Range("C3").Value2 = $readyStatus
$wsOverview.Range("D3").Value2 = $handoffDate
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $png1Url, $Missing, $Missing, $png1) | Out-Null

# Row 4: second new dependency (png)
$wsOverview.Range("B4").Value2 = $readyStatus
$wsOverview.Range("C4").Value2 = $readyStatus
$wsOverview.Range("D4").Value2 = $handoffDate
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $png2Url, $Missing, $Missing, $png2) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 2: refresh file names / dates for the existing handoff
$wsZh.Range("A2").Value2 = $mdName
$wsZh.Range("D2").Value2 = $xlfZh
$wsZh.Range("E2").Value2 = $zhDate
$wsZh.Range("H2").Value2 = $epoch
$wsZh.Range("I2").Value2 = $includeReason

# Row 3: new dependency (png #1)
$wsZh.Range("C3").Value2 = $readyStatus
$wsZh.Range("E3").Value2 = $zhDate
$wsZh.Range("H3").Value2 = $epoch
$wsZh.Range("I3").Value2 = $isDependencyReason
$wsZh.Range("J3").Value2 = $dependencyFrom
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $png1Url, $Missing, $Missing, $png1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), $png1Url, $Missing, $Missing, ".png") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $png1ZhUrl, $Missing, $Missing, $png1TargetZh) | Out-Null

# Row 4: new dependency (png #2)
$wsZh.Range("C4").Value2 = $readyStatus
$wsZh.Range("E4").Value2 = $zhDate
$wsZh.Range("H4").Value2 = $epoch
$wsZh.Range("I4").Value2 = $isDependencyReason
$wsZh.Range("J4").Value2 = $dependencyFrom
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $png2Url, $Missing, $Missing, $png2) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B4"), $png2Url, $Missing, $Missing, ".png") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), $png2ZhUrl, $Missing, $Missing, $png2TargetZh) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 2: refresh file names / dates for the existing handoff
$wsDe.Range("A2").Value2 = $mdName
$wsDe.Range("D2").Value2 = $xlfDe
$wsDe.Range("E2").Value2 = $deDate
$wsDe.Range("H2").Value2 = $epoch
$wsDe.Range("I2").Value2 = $includeReason

# Row 3: new dependency (png #1)
$wsDe.Range("C3").Value2 = $readyStatus
$wsDe.Range("E3").Value2 = $deDate
$wsDe.Range("H3").Value2 = $epoch
$wsDe.Range("I3").Value2 = $isDependencyReason
$wsDe.Range("J3").Value2 = $dependencyFrom
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $png1Url, $Missing, $Missing, $png1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), $png1Url, $Missing, $Missing, ".png") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $png1DeUrl, $Missing, $Missing, $png1TargetZh) | Out-Null

# Row 4: new dependency (png #2)
$wsDe.Range("C4").Value2 = $readyStatus
$wsDe.Range("E4").Value2 = $deDate
$wsDe.Range("H4").Value2 = $epoch
$wsDe.Range("I4").Value2 = $isDependencyReason
$wsDe.Range("J4").Value2 = $dependencyFrom
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $png2Url, $Missing, $Missing, $png2) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B4"), $png2Url, $Missing, $Missing, ".png") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), $png2DeUrl, $Missing, $Missing, $png2TargetZh) | Out-Null

Write-Host "Handoff report generated"
